$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels with units
$ws.Range("J1").Value = 'MAE [$COP/kWh]'
$ws.Range("K1").Value = 'MSE [$COP/kWh]'
$ws.Range("L1").Value = 'RMSE [$COP/kWh]'
$ws.Range("M1").Value = 'MAPE [%]'

# Update row 2 data values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = 25
$ws.Range("G2").Value = "<keras.src.optimizers.adam.Adam object at 0x0000012D02197B50>"
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 69.65316057194121
$ws.Range("K2").Value = 5636.96271815866
$ws.Range("L2").Value = 75.07970909745629
$ws.Range("M2").Value = 42.16618722728619
